$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.944.37'
$ws.Range('D3').Value = '2.362.89'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('E9').Value = '  -2.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.13'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.93%  '
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('E13').Value = '  -3.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.73'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').Value = '2.728.29'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').Value = '2.344.03'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.793'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '42.932.99'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.89'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.01%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').Value = '0.0₃0885'
$ws.Range('E21').Value = '  -0.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.97'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '  -3.43%  '
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('E28').Value = '  +15.18%  '
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.20'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0715'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('E37').Value = '  -11.39%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.85'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.22%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.31'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.80%  '
$ws.Range('E40').Value = '  -2.01%  '
$ws.Range('E41').Value = '  -1.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('D43').Value = '1.931.28'
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('E44').Value = '  -0.74%  '
$ws.Range('E45').Value = '  +4.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.43%  '
$ws.Range('E47').Value = '  -2.37%  '
$ws.Range('D48').Value = '2.589.54'
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('E49').Value = '  +1.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.47'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.61%  '
$ws.Range('E51').Value = '  +0.64%  '
